$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 649, pushing the existing rows 649-703 down to 650-704.
$ws.Rows(649).Insert()

# The row that used to be 649 is now at 650 (with all its original data and
# formatting intact). Duplicate it back into the freshly inserted row 649 so
# the new week's record starts from the same template as the prior entry.
$src = $ws.Range("A650:T650")
$dst = $ws.Range("A649:T649")
$src.Copy($dst)

# Update the date for the newly inserted record (row 649) to the new week.
$ws.Cells.Item(649, 4).Value2 = 45166
